$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the formatting of the other headers
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New era/save data values for H2:H7
$values = @(0, 0, 1, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
